# Update the cached date placeholder text ("10/13/16" -> "10/14/16") everywhere
# it appears: notes master, slide master, and every slide layout.
function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -ne 0) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "10/13/16") {
                $tr.Text = "10/14/16"
            }
        }
    }
}

$p = $ppt.ActivePresentation

# Notes master
Update-DateShape $p.NotesMaster

# Slide master
$master = $p.SlideMaster
Update-DateShape $master

# Every slide layout belonging to the slide master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li)
}

# ---------------------------------------------------------------
# Slide 3: "deletePerson(p)" -> "delete" + "Task" + "(" + "task" + ")"
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(16)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Characters(7,6).Text = "Task"
$tr3.Characters(12,1).Text = "task"

# ---------------------------------------------------------------
# Slide 4: "Task" + "ListPanel" -> "TaskListPanel" (merge to one run)
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(11)
$tr4 = $sh4.TextFrame.TextRange
$tr4.Text = "___placeholder___"
$tr4.Text = "TaskListPanel"

# ---------------------------------------------------------------
# Slide 5: "Find" + "Command" -> "FindCommand" (merge to one run)
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(9)
$tr5 = $sh5.TextFrame.TextRange
$tr5.Text = "___placeholder___"
$tr5.Text = "FindCommand"

# ---------------------------------------------------------------
# Slide 6: "(p" + ")" -> "(p)" (merge to one run, keep "deleteTask" run intact)
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(27)
$tr6 = $sh6.TextFrame.TextRange
$tr6.Characters(11,3).Text = "(p)"
